# Implement vaccinations in the initial conditions and later (#30)
#
# This script:
#  1) Removes the old "helper table" (I169:J175) that was used to convert a
#     measured chart height into a number for the 2021-W05 block.
#  2) Adds three new weekly data blocks (week 7, 8, 9 of 2021), each preceded
#     by a "Source: ..." separator row, mirroring the existing layout.
#  3) Recreates the helper table (now in I184:J190) for the week-7 block,
#     including the same two formulas used previously.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve the special formatting used by the old helper table so we can
# --- reapply it to the new helper table before the old cells are cleared.
$ws.Cells.Item(174, 9).Copy()
$ws.Cells.Item(189, 9).PasteSpecial(-4122)
$ws.Cells.Item(174, 10).Copy()
$ws.Cells.Item(189, 10).PasteSpecial(-4122)

$ws.Cells.Item(175, 9).Copy()
$ws.Cells.Item(190, 9).PasteSpecial(-4122)
$ws.Cells.Item(175, 10).Copy()
$ws.Cells.Item(190, 10).PasteSpecial(-4122)

# --- Remove the old helper table that lived alongside the 2021-W05 block.
$ws.Range("I169:J175").Clear()

# --- Week 7 (2021) block, preceded by its source separator -----------------
$ws.Cells.Item(183, 2).Value = "Source: 2021-03-09"

$ws.Cells.Item(184, 1).Value = 2021
$ws.Cells.Item(184, 2).Value = 7
$ws.Cells.Item(184, 3).Value = "0-4"
$ws.Cells.Item(184, 4).Value = 7500
$ws.Cells.Item(184, 5).Value = 5.6
$ws.Cells.Item(184, 9).Value = "Länge für Einheit (cm)"
$ws.Cells.Item(184, 10).Value = 16.75

$ws.Cells.Item(185, 1).Value = 2021
$ws.Cells.Item(185, 2).Value = 7
$ws.Cells.Item(185, 3).Value = "5-14"
$ws.Cells.Item(185, 4).Value = 11250
$ws.Cells.Item(185, 5).Value = 9.7
$ws.Cells.Item(185, 9).Value = "#:"
$ws.Cells.Item(185, 10).Value = 10

$ws.Cells.Item(186, 1).Value = 2021
$ws.Cells.Item(186, 2).Value = 7
$ws.Cells.Item(186, 3).Value = "15-34"
$ws.Cells.Item(186, 4).Value = 93750
$ws.Cells.Item(186, 5).Value = 5.6
$ws.Cells.Item(186, 9).Value = "# / cm"
$ws.Cells.Item(186, 10).Formula = "=J185 / J184"

$ws.Cells.Item(187, 1).Value = 2021
$ws.Cells.Item(187, 2).Value = 7
$ws.Cells.Item(187, 3).Value = "35-59"
$ws.Cells.Item(187, 4).Value = 150000
$ws.Cells.Item(187, 5).Value = 5.5
$ws.Cells.Item(187, 9).Value = "Achsenabschnitt (cm)"
$ws.Cells.Item(187, 10).Value = 0

$ws.Cells.Item(188, 1).Value = 2021
$ws.Cells.Item(188, 2).Value = 7
$ws.Cells.Item(188, 3).Value = "60-79"
$ws.Cells.Item(188, 4).Value = 84750
$ws.Cells.Item(188, 5).Value = 5.7
$ws.Cells.Item(188, 9).Value = "Achsenabschnitt (#)"
$ws.Cells.Item(188, 10).Value = 5

$ws.Cells.Item(189, 1).Value = 2021
$ws.Cells.Item(189, 2).Value = 7
$ws.Cells.Item(189, 3).Value = ">=80"
$ws.Cells.Item(189, 4).Value = 46500
$ws.Cells.Item(189, 5).Value = 7.5
$ws.Cells.Item(189, 9).Value = "Gemessene Höhe (cm)"
$ws.Cells.Item(189, 10).Value = 7.45

$ws.Cells.Item(190, 2).Value = "Source: 2021-03-23"
$ws.Cells.Item(190, 9).Value = "Zahl"
$ws.Cells.Item(190, 10).Formula = "=(J189-J187)*J186+J188"

# --- Week 8 (2021) block -----------------------------------------------------
$ws.Cells.Item(191, 1).Value = 2021
$ws.Cells.Item(191, 2).Value = 8
$ws.Cells.Item(191, 3).Value = "0-4"
$ws.Cells.Item(191, 4).Value = 9333
$ws.Cells.Item(191, 5).Value = 6.7

$ws.Cells.Item(192, 1).Value = 2021
$ws.Cells.Item(192, 2).Value = 8
$ws.Cells.Item(192, 3).Value = "5-14"
$ws.Cells.Item(192, 4).Value = 14000
$ws.Cells.Item(192, 5).Value = 9.8

$ws.Cells.Item(193, 1).Value = 2021
$ws.Cells.Item(193, 2).Value = 8
$ws.Cells.Item(193, 3).Value = "15-34"
$ws.Cells.Item(193, 4).Value = 96666
$ws.Cells.Item(193, 5).Value = 6.6

$ws.Cells.Item(194, 1).Value = 2021
$ws.Cells.Item(194, 2).Value = 8
$ws.Cells.Item(194, 3).Value = "35-59"
$ws.Cells.Item(194, 4).Value = 154000
$ws.Cells.Item(194, 5).Value = 5.9

$ws.Cells.Item(195, 1).Value = 2021
$ws.Cells.Item(195, 2).Value = 8
$ws.Cells.Item(195, 3).Value = "60-79"
$ws.Cells.Item(195, 4).Value = 86667
$ws.Cells.Item(195, 5).Value = 5.7

$ws.Cells.Item(196, 1).Value = 2021
$ws.Cells.Item(196, 2).Value = 8
$ws.Cells.Item(196, 3).Value = ">=80"
$ws.Cells.Item(196, 4).Value = 45000
$ws.Cells.Item(196, 5).Value = 6.6

$ws.Cells.Item(197, 2).Value = "Source: 2021-03-23"

# --- Week 9 (2021) block -----------------------------------------------------
$ws.Cells.Item(198, 1).Value = 2021
$ws.Cells.Item(198, 2).Value = 9
$ws.Cells.Item(198, 3).Value = "0-4"
$ws.Cells.Item(198, 4).Value = 14000
$ws.Cells.Item(198, 5).Value = 6

$ws.Cells.Item(199, 1).Value = 2021
$ws.Cells.Item(199, 2).Value = 9
$ws.Cells.Item(199, 3).Value = "5-14"
$ws.Cells.Item(199, 4).Value = 19333
$ws.Cells.Item(199, 5).Value = 9.4

$ws.Cells.Item(200, 1).Value = 2021
$ws.Cells.Item(200, 2).Value = 9
$ws.Cells.Item(200, 3).Value = "15-34"
$ws.Cells.Item(200, 4).Value = 102666
$ws.Cells.Item(200, 5).Value = 6.3

$ws.Cells.Item(201, 1).Value = 2021
$ws.Cells.Item(201, 2).Value = 9
$ws.Cells.Item(201, 3).Value = "35-59"
$ws.Cells.Item(201, 4).Value = 156666
$ws.Cells.Item(201, 5).Value = 6.1

$ws.Cells.Item(202, 1).Value = 2021
$ws.Cells.Item(202, 2).Value = 9
$ws.Cells.Item(202, 3).Value = "60-79"
$ws.Cells.Item(202, 4).Value = 88666
$ws.Cells.Item(202, 5).Value = 5.3

$ws.Cells.Item(203, 1).Value = 2021
$ws.Cells.Item(203, 2).Value = 9
$ws.Cells.Item(203, 3).Value = ">=80"
$ws.Cells.Item(203, 4).Value = 45333
$ws.Cells.Item(203, 5).Value = 5.3

# --- Update the frozen-pane selection to point at the new bottom of data ---
$ws.Range("E200").Select()
